$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forward/reverse read counts and MD5 checksums for under_ice_rerun samples (rows 12-35)
$ws.Range("BB12").Value = 26132363
$ws.Range("BC12").Value = 26132363
$ws.Range("BD12").Value = "2f43253d68a8a0ffaf0c47d13b4181ea"
$ws.Range("BE12").Value = "8506cce1a9c4bf5ab5a8ef1630bc01ea"
$ws.Range("BB13").Value = 32776170
$ws.Range("BC13").Value = 32776170
$ws.Range("BD13").Value = "9a47e4970377569015ef80b957a9e6c9"
$ws.Range("BE13").Value = "2978194f604c56db0c4c97814da135f7"
$ws.Range("BB14").Value = 15374445
$ws.Range("BC14").Value = 15374445
$ws.Range("BD14").Value = "cabdbe29e207198ffca58832ca1e7c0d"
$ws.Range("BE14").Value = "65d51ae75ab896e1559e8a1ea9937bd1"
$ws.Range("BB15").Value = 20236324
$ws.Range("BC15").Value = 20236324
$ws.Range("BD15").Value = "e9ec4fa534e9d05a59424553bd4ffdf6"
$ws.Range("BE15").Value = "e51c69ea4c87be44b0942490ba2bd643"
$ws.Range("BB16").Value = 23475092
$ws.Range("BC16").Value = 23475092
$ws.Range("BD16").Value = "d1b24965eb3cfd018845da67fa1850bc"
$ws.Range("BE16").Value = "f28a721fe50bbe88c90a76869aa6b1c4"
$ws.Range("BB17").Value = 19019111
$ws.Range("BC17").Value = 19019111
$ws.Range("BD17").Value = "ac866f11b4d1ef3083d8769da6203a67"
$ws.Range("BE17").Value = "4d4bbd575ccbfb50261abc8576378caa"
$ws.Range("BB18").Value = 23417974
$ws.Range("BC18").Value = 23417974
$ws.Range("BD18").Value = "df0446869000b3a6b7b82a421a8e5273"
$ws.Range("BE18").Value = "7916a03e672fbe8e9c6551952c37a702"
$ws.Range("BB19").Value = 21959015
$ws.Range("BC19").Value = 21959015
$ws.Range("BD19").Value = "1c078aa61a58baf35d57ceb5366e52fc"
$ws.Range("BE19").Value = "fa5430f711b88ec1de82a3263faea716"
$ws.Range("BB20").Value = 19349699
$ws.Range("BC20").Value = 19349699
$ws.Range("BD20").Value = "3b323211e3f52500005ee704d521cff4"
$ws.Range("BE20").Value = "0991a0d7d509cc9bb61eb241fdc98798"
$ws.Range("BB21").Value = 24640452
$ws.Range("BC21").Value = 24640452
$ws.Range("BD21").Value = "0cfc3ee43c7964cbdad2237f0b68db64"
$ws.Range("BE21").Value = "955647fea0b807e75ac16a2b1140adf3"
$ws.Range("BB22").Value = 22292795
$ws.Range("BC22").Value = 22292795
$ws.Range("BD22").Value = "ed6c1bb21e7fd960a37ceb39fbec6933"
$ws.Range("BE22").Value = "d80781d60a135d10479589a434640c08"
$ws.Range("BB23").Value = 24336672
$ws.Range("BC23").Value = 24336672
$ws.Range("BD23").Value = "e2a26db2b8e08bcf7e2f81c5f7e2972c"
$ws.Range("BE23").Value = "9ee99e92c1277e5020faa4abed2d0699"
$ws.Range("BB24").Value = 22153770
$ws.Range("BC24").Value = 22153770
$ws.Range("BD24").Value = "49f4330718aba000f342fc250a3b4dbc"
$ws.Range("BE24").Value = "30f67c479152c841478b53718b9cb305"
$ws.Range("BB25").Value = 23514431
$ws.Range("BC25").Value = 23514431
$ws.Range("BD25").Value = "41c4951dac2c197bd4d5ed07ff59914e"
$ws.Range("BE25").Value = "c5b5921d79a97b670562f47fe395677d"
$ws.Range("BB26").Value = 25239007
$ws.Range("BC26").Value = 25239007
$ws.Range("BD26").Value = "f54b725a95775ff92b9f3d1e06829558"
$ws.Range("BE26").Value = "34b7bf3eef1f958258deb92846cfcbb3"
$ws.Range("BB27").Value = 26888537
$ws.Range("BC27").Value = 26888537
$ws.Range("BD27").Value = "2bd49f0fdfe19878074e39e8c86b9789"
$ws.Range("BE27").Value = "3931552e9f9f68a7bd42a87aa6538ce3"
$ws.Range("BB28").Value = 24099741
$ws.Range("BC28").Value = 24099741
$ws.Range("BD28").Value = "318171f4834cf06049826225c22b9aa2"
$ws.Range("BE28").Value = "c988129eec03c2cc4d047bb610347f24"
$ws.Range("BB29").Value = 25131709
$ws.Range("BC29").Value = 25131709
$ws.Range("BD29").Value = "a9d2dc7e4792d06420644734e4b42630"
$ws.Range("BE29").Value = "e30a85a3f83d4ee1b87d16fb23897556"
$ws.Range("BB30").Value = 25333059
$ws.Range("BC30").Value = 25333059
$ws.Range("BD30").Value = "82a2d3b363ed5c78101fdd96d8cd360a"
$ws.Range("BE30").Value = "592f0b9322373b09995f66291928b1f0"
$ws.Range("BB31").Value = 23643065
$ws.Range("BC31").Value = 23643065
$ws.Range("BD31").Value = "ed866a931f11e7d19540321eab158507"
$ws.Range("BE31").Value = "7febe95bea2f13b95a0bfa370c692d0c"
$ws.Range("BB32").Value = 23126069
$ws.Range("BC32").Value = 23126069
$ws.Range("BD32").Value = "a0cf5996f52d9fd7f9ef37a04a725196"
$ws.Range("BE32").Value = "a8bbf5740c6e14d1e77e62a16b44e135"
$ws.Range("BB34").Value = 21665592
$ws.Range("BC34").Value = 21665592
$ws.Range("BD34").Value = "f471246e7f946d08bd763fa512f56cd0"
$ws.Range("BE34").Value = "f53576a0f5f399716c5c8514371b0d59"
$ws.Range("BB35").Value = 24804911
$ws.Range("BC35").Value = 24804911
$ws.Range("BD35").Value = "b913c39657ba36175b85b4e43343eb9b"
$ws.Range("BE35").Value = "a7057c9bdb2f5bbc70b817dcd445f5bb"
